# Update handback/handoff timestamps for the c0a03d35... row (row 4) on both
# the zh-cn and de-de worksheets, as part of generating the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-22 12:50:02"
$wsZhCn.Range("H4").Value = "2016-03-22 12:50:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-22 12:50:11"
$wsDeDe.Range("H4").Value = "2016-03-22 12:50:47"
